$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the "Código" column first for both new rows (matches shared-string
# insertion order observed in the source file), then the remaining columns.
$ws.Range("D5").Value = "uapmarcelino"
$ws.Range("D6").Value = "uapburman"

$ws.Range("A5").Value = "UAP Marcelino"
$ws.Range("B5").Value = "General Lavalle 1583"
$ws.Range("C5").Value = "FLORIDA"

$ws.Range("A6").Value = "URI Burman"
$ws.Range("B6").Value = "Ituzaingó 5725"
$ws.Range("C6").Value = "CARAPACHAY"

$ws.Range("E5").Value = -34.532160508112803
$ws.Range("F5").Value = -58.480465826746901

$ws.Range("E6").Value = -34.527681854771501
$ws.Range("F6").Value = -58.536327416546001

# Match the style of column A cells (vertical-top alignment), same as A1/A2/A4
$ws.Range("A5").VerticalAlignment = -4160 # xlTop
$ws.Range("A6").VerticalAlignment = -4160 # xlTop

# Update selection to match final state
$ws.Range("E7").Select()
